$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.440.14'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '1.699.32'
$ws.Range("E3").Value = '  +0.91%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5464'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.20%  '
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2743'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06451'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.04'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07700'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.98%  '
$ws.Range("D12").Value = '1.694.86'
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.553'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5838'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008393'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.42%  '
$ws.Range("D17").Value = '26.494.01'
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.951'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.257'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.012'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.32%  '
$ws.Range("E25").Value = '  +7.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.907'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06244'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.380'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.54%  '
$ws.Range("E30").Value = '  +0.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.616'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.604'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.691'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.041'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6183'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.416'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.773'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01644'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.43%  '
$ws.Range("D39").Value = '1.119.37'
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("E40").Value = '  -3.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8807'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("E43").Value = '  +0.84%  '
$ws.Range("D44").Value = '1.853.26'
$ws.Range("E44").Value = '  +1.16%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.63'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000108'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.256'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05291'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.155'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.24%  '
